$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns for each changed coin row.
# Column D values that look numeric (e.g. "512.40") are protected with a
# temporary Text number format so Excel keeps them as text (matching the
# original inlineStr cell type) instead of silently converting them to
# numbers; the style is then restored to Normal so no stray formatting is left behind.

$ws.Range("D2").Value = "56.468.96"
$ws.Range("E2").Value = "  +0.01%  "
$ws.Range("D3").Value = "2.328.62"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  -0.30%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "512.40"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "131.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.79%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.35%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.533"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.45%  "
$ws.Range("E9").Value = "  -3.00%  "
$ws.Range("E10").Value = "  -0.58%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.29"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.73%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.339"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").Value = "2.745.59"
$ws.Range("E13").Value = "  -0.22%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.52"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "56.442.51"
$ws.Range("E15").Value = "  -0.17%  "
$ws.Range("E16").Value = "  -1.26%  "
$ws.Range("D17").Value = "2.330.46"
$ws.Range("E17").Value = "  -0.68%  "
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "324.78"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.81%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.57%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("E22").Value = "  -0.06%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "61.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.21%  "
$ws.Range("E24").Value = "  +11.80%  "
$ws.Range("E25").Value = "  +1.50%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.94%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "168.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.26%  "
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("D30").Value = "0.0₃0719"
$ws.Range("E30").Value = "  -3.13%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.31"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.33%  "
$ws.Range("E33").Value = "  +0.03%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.46%  "
$ws.Range("E35").Value = "  +1.46%  "
$ws.Range("E36").Value = "  -1.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.887"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.75%  "
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("E39").Value = "  +1.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "151.25"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +10.68%  "
$ws.Range("E41").Value = "  -0.45%  "
$ws.Range("E42").Value = "  -0.39%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "278.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.77%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0926"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("E46").Value = "  -1.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.558"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.98%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "17.99"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.98%  "
$ws.Range("E49").Value = "  -1.56%  "
$ws.Range("E50").Value = "  +0.07%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.47%  "
